# Add a new condition row (row 25) to the "condition.csv" sheet,
# introducing three new shared strings: itemCanEquip / 道具可以装备 / category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "itemCanEquip"
$ws.Range("B25").Value = "道具可以装备"
$ws.Range("C25").Value = "item"
$ws.Range("D25").Value = "category"
$ws.Range("E25").Value = "<="
$ws.Range("F25").Value = "number"
$ws.Range("G25").Value = 3

$ws.Range("E25").Select()
